# Proyecto 1 - Presupuesto
# Update the "DPI" (A2) and "Nacimiento" (B2) values on Hoja1, move the
# active selection to C6, and switch the sheet's page orientation to
# portrait (this is what stamps a <pageSetup> record onto the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: DPI number (note the trailing tab that was typed/pasted by the user)
$ws.Range("A2").Value = "3306466721202`t"

# B2: Fecha de nacimiento
$ws.Range("B2").Value = "18/09/1997"

# Configure the page as Portrait, which is what Excel records in
# <pageSetup orientation="portrait" .../> when the sheet is saved.
$ws.PageSetup.Orientation = 1

# Leave the cursor on C6, matching the saved cursor/selection position.
$ws.Range("C6").Select()
